# Extend the Binance ETHUSDT OHLCV sheet with 5 more daily rows (16-20)
# and update the existing row 15 (index 13) with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Force the "numeric-looking" text columns to Text format before
#    typing into them, otherwise Excel auto-converts strings such as
#    "126.00000000" into the number 126 and the trailing zeros / decimal
#    formatting captured by the source data is lost.
#    Columns: B(open) C(high) D(low) E(close) F(volume) H(quote vol)
#             J(taker buy base) K(taker buy quote) L(open time GST)
#    (G = Close_time and I = number of trades stay true numbers.)
#    NOTE: multi-area ("A1:A2,B1:B2") ranges don't reliably propagate
#    NumberFormat writes here, so each contiguous block is set alone.
# ---------------------------------------------------------------------
$textRange1 = $ws.Range("B15:F20")
$textRange2 = $ws.Range("H15:H20")
$textRange3 = $ws.Range("J15:L20")
$textRange1.NumberFormat = "@"
$textRange2.NumberFormat = "@"
$textRange3.NumberFormat = "@"

# ---------------------------------------------------------------------
# 2) Row 15 (data index 13) - corrected values
# ---------------------------------------------------------------------
$ws.Range("D15").Value = "126.00000000"
$ws.Range("E15").Value = "127.80000000"
$ws.Range("F15").Value = "421600.75655000"
$ws.Range("H15").Value = "55501305.92012300"
$ws.Range("I15").Value = 145867
$ws.Range("J15").Value = "194731.21180000"
$ws.Range("K15").Value = "25638865.78568890"
$ws.Range("M15").Value = 127.7999999999999
$ws.Range("N15").Value = 129.9450000000001
$ws.Range("O15").Value = 128.2614285714286
$ws.Range("R15").Value = 127.8
$ws.Range("S15").Value = 128.7174304615879
$ws.Range("T15").Value = 131.5394353129163
$ws.Range("U15").Value = 133.3233479202179
$ws.Range("V15").Value = -1.783912607301687
$ws.Range("W15").Value = -1.464915537082233
$ws.Range("X15").Value = -0.3189970702194544

# ---------------------------------------------------------------------
# 3) Row 16 (data index 14) - new row
# ---------------------------------------------------------------------
$ws.Range("B16").Value = "127.80000000"
$ws.Range("C16").Value = "129.69000000"
$ws.Range("D16").Value = "126.61000000"
$ws.Range("E16").Value = "127.75000000"
$ws.Range("F16").Value = "200637.10098000"
$ws.Range("G16").Value = 1577231999999
$ws.Range("H16").Value = "25667665.14201030"
$ws.Range("I16").Value = 87657
$ws.Range("J16").Value = "100608.48356000"
$ws.Range("K16").Value = "12872420.08887930"
$ws.Range("L16").Value = "2019-12-24 08:00:00"
$ws.Range("M16").Value = 127.7499999999999
$ws.Range("N16").Value = 127.775
$ws.Range("O16").Value = 129.1
$ws.Range("P16").Value = 134.7453333333333
$ws.Range("R16").Value = 127.75
$ws.Range("S16").Value = 128.0724767755814
$ws.Range("T16").Value = 130.9046400251341
$ws.Range("U16").Value = 132.720448196802
$ws.Range("V16").Value = -1.815808171667896
$ws.Range("W16").Value = -1.537653296387676
$ws.Range("X16").Value = -0.2781548752802205

# ---------------------------------------------------------------------
# 4) Row 17 (data index 15) - new row
# ---------------------------------------------------------------------
$ws.Range("B17").Value = "127.70000000"
$ws.Range("C17").Value = "127.84000000"
$ws.Range("D17").Value = "123.07000000"
$ws.Range("E17").Value = "125.09000000"
$ws.Range("F17").Value = "225004.49090000"
$ws.Range("G17").Value = 1577318399999
$ws.Range("H17").Value = "28144469.16072900"
$ws.Range("I17").Value = 93397
$ws.Range("J17").Value = "98563.14610000"
$ws.Range("K17").Value = "12331982.19518460"
$ws.Range("L17").Value = "2019-12-25 08:00:00"
$ws.Range("M17").Value = 125.0899999999999
$ws.Range("N17").Value = 126.42
$ws.Range("O17").Value = 128.0014285714286
$ws.Range("P17").Value = 133.3806666666667
$ws.Range("R17").Value = 125.09
$ws.Range("S17").Value = 126.084158879004
$ws.Range("T17").Value = 129.9437249151322
$ws.Range("U17").Value = 131.9222406327808
$ws.Range("V17").Value = -1.978515717648605
$ws.Range("W17").Value = -1.628379496135982
$ws.Range("X17").Value = -0.3501362215126234

# ---------------------------------------------------------------------
# 5) Row 18 (data index 16) - new row
# ---------------------------------------------------------------------
$ws.Range("B18").Value = "125.09000000"
$ws.Range("C18").Value = "132.26000000"
$ws.Range("D18").Value = "124.32000000"
$ws.Range("E18").Value = "125.58000000"
$ws.Range("F18").Value = "274986.52097000"
$ws.Range("G18").Value = 1577404799999
$ws.Range("H18").Value = "34943780.90920390"
$ws.Range("I18").Value = 103034
$ws.Range("J18").Value = "128410.01213000"
$ws.Range("K18").Value = "16332907.92974440"
$ws.Range("L18").Value = "2019-12-26 08:00:00"
$ws.Range("M18").Value = 125.5799999999999
$ws.Range("N18").Value = 125.335
$ws.Range("O18").Value = 127.6414285714286
$ws.Range("P18").Value = 132.1933333333333
$ws.Range("R18").Value = 125.58
$ws.Range("S18").Value = 125.7480529570654
$ws.Range("T18").Value = 129.2307215690578
$ws.Range("U18").Value = 131.2784478225959
$ws.Range("V18").Value = -2.047726253538031
$ws.Range("W18").Value = -1.714180924009384
$ws.Range("X18").Value = -0.3335453295286464

# ---------------------------------------------------------------------
# 6) Row 19 (data index 17) - new row
# ---------------------------------------------------------------------
$ws.Range("B19").Value = "125.58000000"
$ws.Range("C19").Value = "127.10000000"
$ws.Range("D19").Value = "121.91000000"
$ws.Range("E19").Value = "126.29000000"
$ws.Range("F19").Value = "240012.37451000"
$ws.Range("G19").Value = 1577491199999
$ws.Range("H19").Value = "30020098.88913110"
$ws.Range("I19").Value = 102142
$ws.Range("J19").Value = "115839.00216000"
$ws.Range("K19").Value = "14500702.34866460"
$ws.Range("L19").Value = "2019-12-27 08:00:00"
$ws.Range("M19").Value = 126.2899999999999
$ws.Range("N19").Value = 125.935
$ws.Range("O19").Value = 127.37
$ws.Range("P19").Value = 130.9546666666667
$ws.Range("R19").Value = 126.29
$ws.Range("S19").Value = 126.109350986621
$ws.Range("T19").Value = 128.7547713720667
$ws.Range("U19").Value = 130.7855979726117
$ws.Range("V19").Value = -2.030826600544998
$ws.Range("W19").Value = -1.778671824090809
$ws.Range("X19").Value = -0.2521547764541889

# ---------------------------------------------------------------------
# 7) Row 20 (data index 18) - new row
# ---------------------------------------------------------------------
$ws.Range("B20").Value = "126.28000000"
$ws.Range("C20").Value = "128.59000000"
$ws.Range("D20").Value = "125.84000000"
$ws.Range("E20").Value = "127.28000000"
$ws.Range("F20").Value = "77413.69234000"
$ws.Range("G20").Value = 1577577599999
$ws.Range("H20").Value = "9862089.16080160"
$ws.Range("I20").Value = 34300
$ws.Range("J20").Value = "39247.83833000"
$ws.Range("K20").Value = "5001388.97613230"
$ws.Range("L20").Value = "2019-12-28 08:00:00"
$ws.Range("M20").Value = 127.2799999999999
$ws.Range("N20").Value = 126.7850000000001
$ws.Range("O20").Value = 127.4114285714286
$ws.Range("P20").Value = 129.7866666666667
$ws.Range("R20").Value = 127.28
$ws.Range("S20").Value = 126.8897836628785
$ws.Range("T20").Value = 128.5179772421101
$ws.Range("U20").Value = 130.4476076024285
$ws.Range("V20").Value = -1.929630360318413
$ws.Range("W20").Value = -1.809305001955259
$ws.Range("X20").Value = -0.1203253583631536

# ---------------------------------------------------------------------
# 8) Drop the forced Text number-format back to the sheet default so the
#    cells end up styled exactly like their neighbours (no left-over
#    explicit style index), then restamp the bold/border/centered index
#    header style onto the new A-column index cells (A16:A20), matching
#    the look of A2:A15.
# ---------------------------------------------------------------------
$textRange1.Style = "Normal"
$textRange2.Style = "Normal"
$textRange3.Style = "Normal"

$ws.Range("A15").Copy()
$ws.Range("A16:A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17
$ws.Range("A20").Value = 18
